$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "39.834.20"
$ws.Range("E2").Value = "  +0.06%  "

# Row 3
$ws.Range("D3").Value = "2.208.47"
$ws.Range("E3").Value = "  -0.84%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "288.29"
$ws.Range("E5").Value = "  -2.06%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.17"
$ws.Range("E6").Value = "  +3.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("E7").Value = "  -0.45%  "

# Row 8
$ws.Range("E8").Value = "  -0.16%  "

# Row 9
$ws.Range("E9").Value = "  +0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.34"
$ws.Range("E10").Value = "  +1.23%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0776"
$ws.Range("E11").Value = "  -1.61%  "

# Row 12
$ws.Range("E12").Value = "  +2.27%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.43"
$ws.Range("E13").Value = "  +1.59%  "

# Row 14
$ws.Range("D14").Value = "2.548.87"
$ws.Range("E14").Value = "  -0.93%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.89"
$ws.Range("E15").Value = "  -1.89%  "

# Row 16
$ws.Range("D16").Value = "2.217.15"
$ws.Range("E16").Value = "  -0.39%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.725"
$ws.Range("E17").Value = "  +0.25%  "

# Row 18
$ws.Range("D18").Value = "39.777.83"
$ws.Range("E18").Value = "  +0.05%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.54"
$ws.Range("E19").Value = "  +9.47%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0880"
$ws.Range("E20").Value = "  -0.61%  "

# Row 21
$ws.Range("E21").Value = "  -0.06%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.34"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.15"
$ws.Range("E23").Value = "  +1.10%  "

# Row 24
$ws.Range("E24").Value = "  -0.13%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("E25").Value = "  +0.79%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.82"
$ws.Range("E26").Value = "  -0.67%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.41"
$ws.Range("E27").Value = "  -2.20%  "

# Row 28
$ws.Range("E28").Value = "  -0.94%  "

# Row 29
$ws.Range("E29").Value = "  -0.35%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "152.76"
$ws.Range("E30").Value = "  +1.18%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.56"
$ws.Range("E31").Value = "  -2.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.17%  "

# Row 33
$ws.Range("E33").Value = "  +1.93%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0716"
$ws.Range("E34").Value = "  +1.53%  "

# Row 35
$ws.Range("E35").Value = "  +0.82%  "

# Row 36
$ws.Range("E36").Value = "  +5.33%  "

# Row 37
$ws.Range("E37").Value = "  -0.17%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "15.68"
$ws.Range("E38").Value = "  -2.69%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0982"
$ws.Range("E39").Value = "  +0.21%  "

# Row 40
$ws.Range("E40").Value = "  +2.19%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.84"
$ws.Range("E41").Value = "  +3.37%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.088.01"
$ws.Range("E42").Value = "  +7.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.15"
$ws.Range("E43").Value = "  -0.36%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.97"
$ws.Range("E44").Value = "  +5.96%  "

# Row 45
$ws.Range("E45").Value = "  -0.84%  "

# Row 46
$ws.Range("E46").Value = "  +6.97%  "

# Row 47
$ws.Range("E47").Value = "  +1.53%  "

# Row 48
$ws.Range("D48").Value = "2.422.02"
$ws.Range("E48").Value = "  -0.76%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "88.43"
$ws.Range("E49").Value = "  -0.79%  "

# Row 50
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.44"
$ws.Range("E50").Value = "  +0.24%  "

# Row 51
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "68.77"
$ws.Range("E51").Value = "  -2.94%  "
